# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (people interested) counts in column F for a
# handful of events across all four sheets (the "全部类型" sheet mirrors
# the same events found in the others, so it needs the same bump).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value  = 1217
$ws.Range("F9").Value  = 1069
$ws.Range("F10").Value = 3137
$ws.Range("F17").Value = 263
$ws.Range("F20").Value = 1267
$ws.Range("F21").Value = 1267
$ws.Range("F22").Value = 202
$ws.Range("F29").Value = 617
$ws.Range("F35").Value = 350
$ws.Range("F37").Value = 5150
$ws.Range("F38").Value = 593
$ws.Range("F39").Value = 332
$ws.Range("F40").Value = 219

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 313
$ws.Range("F18").Value = 52
$ws.Range("F23").Value = 748

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 669
$ws.Range("F6").Value = 463

# 全部类型 (All types - aggregate of the sheets above)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 669
$ws.Range("F10").Value = 1217
$ws.Range("F13").Value = 3137
$ws.Range("F20").Value = 263
$ws.Range("F23").Value = 1267
$ws.Range("F24").Value = 1267
$ws.Range("F27").Value = 313
$ws.Range("F30").Value = 52
$ws.Range("F32").Value = 617
$ws.Range("F38").Value = 350
$ws.Range("F42").Value = 332
$ws.Range("F43").Value = 219
